# Apply the "3d_Secure" payment-details test row into the DataSet sheet.
# This mirrors inserting a new row above the existing "Invalid details" row
# (row 9) on the DataSet sheet, right after the "PaymentDetails" row, and
# filling in the new 3D-Secure test card data (OXO US test cases 101/102).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")

# Insert a new blank row above row 9, shifting everything below down by one.
$ws.Rows.Item(9).Insert()

# Populate the new row with the 3D-Secure payment test data.
$ws.Range("A9").Value = "3d_Secure"

# Card number: force text (leading apostrophe) so it is stored as a shared
# string, matching the existing "PaymentDetails" row's card-number cell.
$ws.Range("Y9").Value = "'4000000000003220"

# Expiry: force text too, then restore the date-look "quoted" number format
# (copied from the existing card's expiry cell) so it renders the same way
# without Excel re-interpreting "06/29" as an actual date.
$ws.Range("Z9").Value = "'06/29"
$ws.Range("Z9").NumberFormat = $ws.Range("Z8").NumberFormat

$ws.Range("AA9").Value = 123
$ws.Range("AK9").Value = "$"
